# Update p-values in the "Fig 5" table according to the new data
# produced by preprocess_indicators and indexes.R.
#
# Dimension ecological    : 0.18 -> 0.22
# Dimension institutional : 0.20 -> 0.01
# Dimension socioeconomic : 0.51 -> 0.32
# Species Cod             : 0.58 -> 0.89
# Species Hake            : 0.53 -> 0.24
#
# Each value lives alone in the 3rd column of the single table. The
# paragraph holding the number starts with an empty run (kept for
# formatting continuity) followed by the run that actually carries the
# text. InsertXML lets us replace exactly that paragraph's contents
# in-place (same run/xml:space structure) instead of letting a plain
# Find/Replace collapse the two runs together.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function Set-PValue($row, $newValue) {
    $cell = $tbl.Cell($row, 3)
    $xml = '<w:p xmlns:w="' + $wNs + '">' +
           '<w:pPr><w:pBdr/><w:spacing w:before="20" w:after="20"/><w:ind w:left="20" w:right="20"/><w:jc w:val="center"/></w:pPr>' +
           '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r>' +
           '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">' + $newValue + '</w:t></w:r>' +
           '</w:p>'
    $cell.Range.InsertXML($xml)
}

Set-PValue 2 "0.22"
Set-PValue 3 "0.01"
Set-PValue 4 "0.32"
Set-PValue 5 "0.89"
Set-PValue 6 "0.24"
